$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "average" row label in column A
$ws.Range("A14").Value = "average"

# Average formulas for each numeric column
$ws.Range("B14").Formula = "=AVERAGE(B2:B13)"
$ws.Range("C14:E14").Formula = "=AVERAGE(C2:C13)"

# Green highlight fill for the whole new row
$ws.Range("A14:E14").Interior.Color = 5296274

# Number format (2 decimals) for the numeric cells of the new row
$ws.Range("B14:E14").NumberFormat = "0.00"

# Update selection to match the newly added row
$ws.Range("B14:E14").Select() | Out-Null

$wb.Save() | Out-Null
